$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.133.02'
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = '1.592.13'
$ws.Range("E3").Value = '  -0.98%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9974'
$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9978'
$ws.Range("E5").Value = '  -0.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '301.61'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3772'
$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3616'
$ws.Range("E8").Value = '  -1.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.80'
$ws.Range("E9").Value = '  +2.96%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.248'
$ws.Range("E10").Value = '  -1.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9973'
$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08068'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.39'
$ws.Range("E13").Value = '  -3.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.531'
$ws.Range("E14").Value = '  -1.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.372'
$ws.Range("E15").Value = '  -1.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001240'

$ws.Range("D17").Value = '1.591.44'
$ws.Range("E17").Value = '  -1.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.90'
$ws.Range("E18").Value = '  +1.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06779'
$ws.Range("E19").Value = '  -0.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.94'
$ws.Range("E20").Value = '  -2.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.455'
$ws.Range("E21").Value = '  -1.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9978'
$ws.Range("E22").Value = '  -0.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.81'
$ws.Range("E23").Value = '  -2.08%  '

$ws.Range("D24").Value = '23.120.03'
$ws.Range("E24").Value = '  -0.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.383'
$ws.Range("E25").Value = '  +1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.924'
$ws.Range("E26").Value = '  +1.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.96'
$ws.Range("E27").Value = '  -0.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '148.86'
$ws.Range("E28").Value = '  -1.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.210'
$ws.Range("E29").Value = '  -1.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.98'
$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.376'
$ws.Range("E31").Value = '  -1.95%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.705'
$ws.Range("E32").Value = '  -2.18%  '

$ws.Range("D33").Value = '1.765.66'
$ws.Range("E33").Value = '  -1.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9591'
$ws.Range("E34").Value = '  -1.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07480'
$ws.Range("E35").Value = '  -2.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.13'
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02690'
$ws.Range("E37").Value = '  -2.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2503'
$ws.Range("E38").Value = '  -2.41%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.105'
$ws.Range("E39").Value = '  -2.92%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08793'
$ws.Range("E40").Value = '  -1.30%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7090'
$ws.Range("E41").Value = '  -1.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.355'
$ws.Range("E42").Value = '  -2.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.24'
$ws.Range("E43").Value = '  -4.65%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.96'
$ws.Range("E44").Value = '  -4.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6497'
$ws.Range("E45").Value = '  -2.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9972'
$ws.Range("E46").Value = '  -0.17%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.286'
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.996'
$ws.Range("E48").Value = '  +0.37%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.51'
$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07902'
$ws.Range("E50").Value = '  -1.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.214'
$ws.Range("E51").Value = '  +2.87%  '
